$d = $word.ActiveDocument

$d.Content.Find.Execute("(NOMENCLATURE: ENUM-STATUS_DR)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(NOMENCLATURE: HubSante.etatDemande)", 2)

$d.Content.Find.Execute("(NOMENCLATURE: CISU-CADRE_CONV)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(NOMENCLATURE: HubSante.cadre)", 2)

$d.Content.Find.Execute("(NOMENCLATURE: CISU-Code_Effet_a_obtenir)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(NOMENCLATURE: HubSante.effet)", 2)

$d.Content.Find.Execute("(NOMENCLATURE: SI-SAMU-DELAI)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(NOMENCLATURE: HubSante.delai)", 2)
